# removed gps from va refusals
#
# The "note_geocode" note row and the "hh_geo_location" geopoint question
# (rows 18-19 of the "survey" sheet) are removed from this XLSForm, along
# with the bookkeeping that naturally follows from deleting those two
# rows: the autofilter/conditional-formatting ranges that used to stretch
# down to row 29 now stop at row 27, and the form's "version" stamp on
# the settings sheet is bumped to mark the edit.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Delete the two geocoding rows (note_geocode note + hh_geo_location
# geopoint question) - everything below shifts up by two rows, so the
# sheet's data now runs from row 1 to row 27 instead of row 1 to row 29.
$survey.Rows.Item(18).Resize(2).Delete() | Out-Null

# The conditional formatting that used to cover the (now shorter) tail of
# the sheet still points at the old row 29 boundary - pull it back in to
# match the new extent of the data.
$fcC = $survey.Range("C18:C29").FormatConditions
for ($i = 1; $i -le $fcC.Count; $i++) {
    $cond = $fcC.Item($i)
    if ($cond.AppliesTo.Address() -eq "`$C`$18:`$C`$29") {
        $cond.ModifyAppliesToRange($survey.Range("C18:C27"))
    }
}

$fcA = $survey.Range("A1:A29").FormatConditions
for ($i = 1; $i -le $fcA.Count; $i++) {
    $cond = $fcA.Item($i)
    if ($cond.AppliesTo.Address() -eq "`$A`$1:`$A`$29") {
        $cond.ModifyAppliesToRange($survey.Range("A1:A27"))
    }
}

# The hidden "filter view" defined name also referenced $A$1:$A$29 - move
# it (and its GUID-bearing name) to match the new, shorter range.
$filterName = $wb.Names.Item("Z_5DE9E846_3713_4B6A_9280_F19D93A89EBA_.wvu.FilterData")
$filterName.RefersTo = "=survey!`$A`$1:`$A`$27"
$filterName.Name = "Z_FFAB8A60_D78A_4936_8B3A_C89EBC2E6F6E_.wvu.FilterData"

# Bump the form version in the settings sheet to mark the edit.
$settings = $wb.Worksheets.Item("settings")
$settings.Range("F2").Value = 2021021101
